$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove existing hyperlinks first (rows will be rewritten fully)
$ws.Hyperlinks.Delete()

$data = @{}
$data[2] = @('【急募】LINEとChatGPT GEMINI連携の簡易質問対応システム開発', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5457662', 423, '🔥GPT,ChatGPT ◆開発,システム開発')
$data[3] = @('【PoC】【急募】防犯カメラ用動画AI解析システム開発の専門家募集', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5457475', 403, '🔥AI,Ai ◆開発,システム開発')
$data[4] = @('【急募】AI医療系請求IOSアプリ開発のエキスパート募集', 'システム開発', '500,000 円 ~ 1,000,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5457710', 385, '🔥AI,Ai ◆開発 ◇アプリ')
$data[5] = @('大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集', 'システム開発', '300,000 円 ~ 500,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5423720', 385, '🔥AI,Ai ◆効率化')
$data[6] = @('EC×AIプロダクト/業務改善リード', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5450024', 338, '🔥AI,Ai ◇業務改善')
$data[7] = @('製造業のR&D支援!「プロセスデータ解析」「音響異常検知」のAIエンジニア募集', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5439165', 303, '🔥AI,Ai')
$data[8] = @('【急募】AWSスクレイピングツールの開発を依頼したいです!', 'システム開発', '5,000 円 ~ 10,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5457255', 165, '◆ツール,開発')
$data[9] = @('【フリーランス募集】Webサービス・業務システム開発エンジニア', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5457382', 118, '◆開発,システム開発')
$data[10] = @('【急募】n8nを使った請求書自動化プロジェクトの依頼', 'システム開発', '200,000 円 ~ 300,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5457696', 88, '◆自動化')
$data[11] = @('GoogleCloudを利用したアジャイル開発共通基盤のSREエンジニアの募集', 'システム開発', '500,000 円 ~ 1,000,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5457458', 75, '◆開発')
$data[12] = @('【急募】飲食店予約サイトの制作と将来的なアプリ化(アプリ化の際は別契約)', 'システム開発', '1,000,000 円 ~ 3,000,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5457089', 70, '◇アプリ')
$data[13] = @('【外国人大歓迎】【急募】ECツールの保守・バグ修正・機能追加エンジニア募集', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5457026', 68, '◆ツール')
$data[14] = @('【急募】PHPによる申請サイト構築支援!', 'システム開発', '50,000 円 ~ 100,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5457023', 58, '◇サイト ○PHP')
$data[15] = @('初回 作成中のwebサイトにて会員登録〜ユーザープロフィール情報の安全な保存までを実装できるエンジニアの方', 'システム開発', '10,000 円 ~ 20,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5457692', 30, '◇サイト')
$data[16] = @('【急募】Kintoneでの請求書自動発行システム構築依頼', 'システム開発', '20,000 円 ~ 50,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5457134', 28, $null)
$data[17] = @('【急募】簡単なHP作成とAWS構築をしてくれる方募集', 'システム開発', '50,000 円 ~ 100,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5457524', 18, $null)
$data[18] = @('回路設計者募集|UVA浄化装置 (マイコン不使用/555タイマー制御) ※成果物全帰属', 'システム開発', '100,000 円 ~ 200,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5457451', 18, $null)
$data[19] = @('【急募】UTAGE構築方法をzoomで教えていただける方', 'システム開発', '1,000 ~ 5,000 円 / 固定', '期限情報なし', 'https://www.lancers.jp/work/detail/5457448', 10, $null)

$newDate = '2025-12-20 01:20:31'

foreach ($r in 2..19) {
    $row = $data[$r]
    $ws.Cells.Item($r, 1).Value = $newDate
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    if ($row[6] -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row[6]
    } else {
        $ws.Cells.Item($r, 8).ClearContents()
    }
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row[4])
    $ws.Cells.Item($r, 6).Style = "Hyperlink"
}

Write-Output "done"